$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

# Column L ("Diferencia Stock") should match column K ("Stock Minimo Objetivo")
# for each data row from 3 to 35.
$total = 0
for ($row = 3; $row -le 35; $row++) {
    $kValue = $ws.Cells.Item($row, 11).Value2
    $ws.Cells.Item($row, 12).Value2 = $kValue
    $total = $total + $kValue
}

# Update the total in C49 to be the sum of the "Diferencia Stock" column (L3:L35)
$ws.Range("C49").Value2 = $total
